$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new word entry "vault" on row 18
$ws.Range("A18").Value = "vault"
$ws.Range("B18").Value = "/vO:lt/"
$ws.Range("C18").Value = "vt. vi."

# Update the active selection to C23 (matches recorded cursor position in diff)
$ws.Range("C23").Select()
